$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C rows 2 through 27 hold a date-serial value (last-changed date)
# that was bumped by one day (45271 -> 45272, i.e. 2023-12-11 -> 2023-12-12).
for ($row = 2; $row -le 27; $row++) {
    $ws.Cells.Item($row, 3).Value = 45272
}
